# Add a new "Scenario5" test-case sheet at the end of the workbook,
# cloning the look-and-feel of the existing "Scenario4" sheet (same
# title / header-row formatting) and filling in the new scenario's
# title + standard column headers.

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Scenario5"

# Template sheet to copy header formatting from.
$template = $wb.Worksheets.Item("Scenario4")

# --- Clone formatting for the title row (A1:J1) and the column-header
#     row (A2:J2) in one shot. ---------------------------------------------
$template.Activate()
$template.Range("A1:J2").Select()
$template.Range("A1:J2").Copy()
$ws5.Range("A1:J2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Content ---------------------------------------------------------------
$ws5.Range("A1").Value = "Password Protection with Information Existing in Projection"

$ws5.Range("A2").Value = "Test Case Title "
$ws5.Range("B2").Value = "Test Case ID"
$ws5.Range("C2").Value = "Test Description"
$ws5.Range("D2").Value = "Preconditions"
$ws5.Range("E2").Value = "Test Steps"
$ws5.Range("F2").Value = "Verify Response"
$ws5.Range("G2").Value = "Expected Results"
$ws5.Range("H2").Value = "Postconditions"
$ws5.Range("I2").Value = "Test Data"
$ws5.Range("J2").Value = "Notes"

# Title spans the whole header band, same as the sibling scenario sheets.
$ws5.Range("A1:J1").Merge($false)

# Header row wraps onto multiple lines, same as the sibling scenario sheets.
$ws5.Rows.Item(2).RowHeight = 43.2

# Leave the workbook on the newly added sheet / cell it was saved on
# upstream.
$ws5.Activate()
$ws5.Range("F7").Select()

Write-Output "Scenario5 sheet added"
